$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9508527517318726
$ws.Range("B1").Value = 2.053400039672852
$ws.Range("C1").Value = 7.779109477996826
$ws.Range("D1").Value = 2.632189989089966
$ws.Range("E1").Value = 0.8011994361877441
